# Season-record columns: the previous export only pulled team statistics,
# not the team's Wins / Losses / Ties for the season. Add those three
# columns (AC:AE) to the player table, matching the header style already
# used by the other column headers in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row currently in the sheet (column A has a value on every row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Headers (row 1) -------------------------------------------------
# Copy the formatting of the existing last header cell (AB1: bold, bordered,
# centered/top-aligned) onto the three new header cells so they pick up the
# same style index instead of minting a new one.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Season record values (rows 2..last) -----------------------------
$dataRange = "2:" + $lastRow
$ws.Range("AC2:AC" + $lastRow).Value = 75
$ws.Range("AD2:AD" + $lastRow).Value = 87
$ws.Range("AE2:AE" + $lastRow).Value = 0
